$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    @([double]"0.9903733673316458", [double]"-0.0002537796048658035", [double]"1.40006462893066", [double]"0.1709853632432805", [double]"1.570796390562869", [double]"-0.5804230367656859"),
    @([double]"1.03574262259447", [double]"-9.196697643037226e-05", [double]"1.398862850774126", [double]"0.1720253258411236", [double]"1.57079639084348", [double]"-0.5350537819232115"),
    @([double]"1.239026965236186", [double]"0.0006330608780510351", [double]"1.393478087873858", [double]"0.1766850477574261", [double]"1.570796392100804", [double]"-0.3317694411649427"),
    @([double]"1.525028734600179", [double]"0.001653106231889857", [double]"1.385902237726871", [double]"0.1832408340786359", [double]"1.57079639386974", [double]"-0.04576767445078167"),
    @([double]"1.728313077241895", [double]"0.002378134086371267", [double]"1.380517474826603", [double]"0.1879005559949384", [double]"1.570796395127064", [double]"0.1575166663074876"),
    @([double]"1.77368233250472", [double]"0.002539946714806697", [double]"1.379315696670069", [double]"0.1889405185927815", [double]"1.570796395407675", [double]"0.2028859211499613")
)

for ($r = 1; $r -le 6; $r++) {
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $values[$r - 1][$c - 1]
    }
}
